$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (order chosen to reproduce the target shared-string table order) ---
$ws.Range("C1").Value = "Categoria"
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "Marca"

# --- Data row ---
$ws.Range("A2").Value = "Mesa de Centro Moderna Montego 1 Cajón"
$ws.Range("B2").Value = "TU MESITA"
$ws.Range("C2").Value = "Sala"

# --- Number formats: make A2:B2 / C2 / B3 use the "Text" (@) style like A2 already had ---
$ws.Range("A2:B2").ClearFormats()
$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("B3").NumberFormat = "@"

# --- Column widths: columns A and B share the same custom width ---
$ws.Range("A1:B1").ColumnWidth = 42

# --- Selection matches the saved view state ---
$null = $ws.Range("B13").Select()
